$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 17071
$ws1.Range("F9").Value = 915
$ws1.Range("F11").Value = 237
$ws1.Range("F13").Value = 11805
$ws1.Range("F14").Value = 32
$ws1.Range("F15").Value = 34
$ws1.Range("F16").Value = 1487
$ws1.Range("F17").Value = 4689
$ws1.Range("F18").Value = 488
$ws1.Range("F25").Value = 41

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 17071
$ws4.Range("F9").Value = 917
$ws4.Range("F11").Value = 237
$ws4.Range("F15").Value = 11805
$ws4.Range("F16").Value = 32
$ws4.Range("F17").Value = 34
$ws4.Range("F18").Value = 1487
$ws4.Range("F19").Value = 4689
$ws4.Range("F20").Value = 488
$ws4.Range("F27").Value = 41
